$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(54, 8).Value = 25357.143
$ws.Cells.Item(54, 10).Value = 28750
$ws.Cells.Item(54, 12).Value = 28750
$ws.Cells.Item(54, 14).Value = -29722
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = 0
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = 0
$ws.Cells.Item(100, 8).Value = 946.4375
$ws.Cells.Item(100, 9).Value = 982.7778
$ws.Cells.Item(100, 10).Value = 899.7143
$ws.Cells.Item(100, 11).Value = 982.7778
$ws.Cells.Item(100, 12).Value = 899.7143
$ws.Cells.Item(100, 13).Value = -441.7778
$ws.Cells.Item(100, 14).Value = -1981.7143
$ws.Cells.Item(132, 8).Value = 1849.4642
$ws.Cells.Item(132, 9).Value = 1828.5416
$ws.Cells.Item(132, 11).Value = 5485.6248
$ws.Cells.Item(132, 13).Value = -2955.6248
$ws.Cells.Item(137, 8).Value = 8848.643
$ws.Cells.Item(137, 9).Value = 8038.467
$ws.Cells.Item(137, 10).Value = 9783.462
$ws.Cells.Item(137, 11).Value = 24115.401
$ws.Cells.Item(137, 12).Value = 29350.386
$ws.Cells.Item(137, 13).Value = -21565.401
$ws.Cells.Item(137, 14).Value = -34450.386
$ws.Cells.Item(141, 8).Value = 2635.923
$ws.Cells.Item(141, 9).Value = 1326.7
$ws.Cells.Item(141, 10).Value = 7000
$ws.Cells.Item(141, 11).Value = 3980.1
$ws.Cells.Item(141, 12).Value = 21000
$ws.Cells.Item(141, 13).Value = 1199.9
$ws.Cells.Item(141, 14).Value = -31360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 32926.184
$ws.Cells.Item(61, 9).Value = 54202.285
$ws.Cells.Item(61, 10).Value = 6643.9414
$ws.Cells.Item(61, 11).Value = 54202.285
$ws.Cells.Item(61, 12).Value = 6643.9414
$ws.Cells.Item(61, 13).Value = -53990.285
$ws.Cells.Item(61, 14).Value = -7067.9414
$ws.Cells.Item(74, 8).Value = 23424.271
$ws.Cells.Item(74, 9).Value = 2038.1471
$ws.Cells.Item(74, 10).Value = 75362
$ws.Cells.Item(74, 11).Value = 2038.1471
$ws.Cells.Item(74, 12).Value = 75362
$ws.Cells.Item(74, 13).Value = -1164.1471
$ws.Cells.Item(74, 14).Value = -77110
$ws.Cells.Item(77, 8).Value = 23424.271
$ws.Cells.Item(77, 9).Value = 2038.1471
$ws.Cells.Item(77, 10).Value = 75362
$ws.Cells.Item(77, 11).Value = 10190.7355
$ws.Cells.Item(77, 12).Value = 376810
$ws.Cells.Item(77, 13).Value = -5822.735499999999
$ws.Cells.Item(77, 14).Value = -385546
$ws.Cells.Item(136, 8).Value = 32926.184
$ws.Cells.Item(136, 9).Value = 54202.285
$ws.Cells.Item(136, 10).Value = 6643.9414
$ws.Cells.Item(136, 11).Value = 162606.855
$ws.Cells.Item(136, 12).Value = 19931.8242
$ws.Cells.Item(136, 13).Value = -160056.855
$ws.Cells.Item(136, 14).Value = -25031.8242

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 4321.558
$ws.Cells.Item(86, 10).Value = 5062.24
$ws.Cells.Item(86, 12).Value = 5062.24
$ws.Cells.Item(86, 14).Value = -7308.24
$ws.Cells.Item(89, 8).Value = 4321.558
$ws.Cells.Item(89, 10).Value = 5062.24
$ws.Cells.Item(89, 12).Value = 25311.2
$ws.Cells.Item(89, 14).Value = -36543.2
$ws.Cells.Item(95, 8).Value = 55000
$ws.Cells.Item(95, 10).Value = 55000
$ws.Cells.Item(95, 12).Value = 55000
$ws.Cells.Item(95, 14).Value = -60492
$ws.Cells.Item(107, 8).Value = 42670.168
$ws.Cells.Item(107, 9).Value = 53755.5
$ws.Cells.Item(107, 10).Value = 20499.5
$ws.Cells.Item(107, 11).Value = 53755.5
$ws.Cells.Item(107, 12).Value = 20499.5
$ws.Cells.Item(107, 13).Value = -51835.5
$ws.Cells.Item(107, 14).Value = -24339.5
$ws.Cells.Item(134, 8).Value = 1839.2941
$ws.Cells.Item(134, 9).Value = 1946.85
$ws.Cells.Item(134, 10).Value = 1685.6428
$ws.Cells.Item(134, 11).Value = 5840.549999999999
$ws.Cells.Item(134, 12).Value = 5056.928400000001
$ws.Cells.Item(134, 13).Value = -3305.549999999999
$ws.Cells.Item(134, 14).Value = -10126.9284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1531.8334
$ws.Cells.Item(58, 9).Value = 1327.125
$ws.Cells.Item(58, 11).Value = 1327.125
$ws.Cells.Item(58, 13).Value = -1124.125
$ws.Cells.Item(107, 8).Value = 2139.5386
$ws.Cells.Item(107, 9).Value = 2376
$ws.Cells.Item(107, 11).Value = 2376
$ws.Cells.Item(107, 13).Value = -456
$ws.Cells.Item(136, 8).Value = 1531.8334
$ws.Cells.Item(136, 9).Value = 1327.125
$ws.Cells.Item(136, 11).Value = 3981.375
$ws.Cells.Item(136, 13).Value = -1431.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 622.6875
$ws.Cells.Item(5, 9).Value = 622.6875
$ws.Cells.Item(5, 11).Value = 1868.0625
$ws.Cells.Item(5, 13).Value = -1756.0625
$ws.Cells.Item(105, 8).Value = 14993
$ws.Cells.Item(105, 10).Value = 16329.556
$ws.Cells.Item(105, 12).Value = 48988.66800000001
$ws.Cells.Item(105, 14).Value = -54230.66800000001
$ws.Cells.Item(131, 8).Value = 26542.596
$ws.Cells.Item(131, 9).Value = 1577.3334
$ws.Cells.Item(131, 10).Value = 33351.305
$ws.Cells.Item(131, 11).Value = 4732.0002
$ws.Cells.Item(131, 12).Value = 100053.915
$ws.Cells.Item(131, 13).Value = 307.9997999999996
$ws.Cells.Item(131, 14).Value = -110133.915
$ws.Cells.Item(132, 8).Value = 2250
$ws.Cells.Item(132, 10).Value = 2400
$ws.Cells.Item(132, 12).Value = 21600
$ws.Cells.Item(132, 14).Value = -26660
$ws.Cells.Item(135, 8).Value = 622.6875
$ws.Cells.Item(135, 9).Value = 622.6875
$ws.Cells.Item(135, 11).Value = 5604.1875
$ws.Cells.Item(135, 13).Value = -3069.1875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1532.7037
$ws.Cells.Item(102, 9).Value = 1463.619
$ws.Cells.Item(102, 11).Value = 1463.619
$ws.Cells.Item(102, 13).Value = 158.3810000000001
$ws.Cells.Item(123, 8).Value = 40249.25
$ws.Cells.Item(123, 9).Value = 37999
$ws.Cells.Item(123, 10).Value = 40999.332
$ws.Cells.Item(123, 11).Value = 37999
$ws.Cells.Item(123, 12).Value = 40999.332
$ws.Cells.Item(123, 13).Value = -35549
$ws.Cells.Item(123, 14).Value = -45899.332
$ws.Cells.Item(134, 8).Value = 84883.664
$ws.Cells.Item(134, 10).Value = 84883.664
$ws.Cells.Item(134, 12).Value = 254650.992
$ws.Cells.Item(134, 14).Value = -259720.992
$ws.Cells.Item(141, 8).Value = 85122
$ws.Cells.Item(141, 10).Value = 85122
$ws.Cells.Item(141, 12).Value = 85122
$ws.Cells.Item(141, 14).Value = -95482

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1019.4
$ws.Cells.Item(16, 9).Value = 699
$ws.Cells.Item(16, 10).Value = 1500
$ws.Cells.Item(16, 11).Value = 699
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = -529
$ws.Cells.Item(16, 14).Value = -1840
$ws.Cells.Item(63, 8).Value = 46598
$ws.Cells.Item(66, 8).Value = 46598
$ws.Cells.Item(82, 8).Value = 604
$ws.Cells.Item(82, 9).Value = 665.8570999999999
$ws.Cells.Item(82, 10).Value = 495.75
$ws.Cells.Item(82, 11).Value = 665.8570999999999
$ws.Cells.Item(82, 12).Value = 495.75
$ws.Cells.Item(82, 13).Value = -304.8570999999999
$ws.Cells.Item(82, 14).Value = -1217.75
$ws.Cells.Item(85, 8).Value = 604
$ws.Cells.Item(85, 9).Value = 665.8570999999999
$ws.Cells.Item(85, 10).Value = 495.75
$ws.Cells.Item(85, 11).Value = 665.8570999999999
$ws.Cells.Item(85, 12).Value = 495.75
$ws.Cells.Item(85, 13).Value = 582.1429000000001
$ws.Cells.Item(85, 14).Value = -2991.75
$ws.Cells.Item(100, 8).Value = 3166.1
$ws.Cells.Item(100, 9).Value = 3166.1
$ws.Cells.Item(100, 11).Value = 3166.1
$ws.Cells.Item(100, 13).Value = -2625.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 7909.75
$ws.Cells.Item(41, 10).Value = 8297.546
$ws.Cells.Item(41, 12).Value = 8297.546
$ws.Cells.Item(41, 14).Value = -9077.546
$ws.Cells.Item(107, 8).Value = 1192891.6
$ws.Cells.Item(107, 10).Value = 2043341.9
$ws.Cells.Item(107, 12).Value = 6130025.699999999
$ws.Cells.Item(107, 14).Value = -6133865.699999999
$ws.Cells.Item(126, 8).Value = 19235614
$ws.Cells.Item(126, 10).Value = 4297.8
$ws.Cells.Item(126, 12).Value = 12893.4
$ws.Cells.Item(126, 14).Value = -17833.4
$ws.Cells.Item(132, 8).Value = 2413.5908
$ws.Cells.Item(132, 9).Value = 2137.2334
$ws.Cells.Item(132, 10).Value = 3005.7856
$ws.Cells.Item(132, 11).Value = 6411.7002
$ws.Cells.Item(132, 12).Value = 9017.356800000001
$ws.Cells.Item(132, 13).Value = -3881.7002
$ws.Cells.Item(132, 14).Value = -14077.3568
$ws.Cells.Item(140, 8).Value = 77196.39999999999
$ws.Cells.Item(140, 10).Value = 77196.39999999999
$ws.Cells.Item(140, 12).Value = 77196.39999999999
$ws.Cells.Item(140, 14).Value = -87556.39999999999
